$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lookup table: "scope/<X>" entries become "Additional Information/<X>"
$ws.Range("A14").Value = "view/metadata/custom_fields/Additional Information/Economy Coverage"
$ws.Range("A20").Value = "view/metadata/custom_fields/Additional Information/Update Frequency"
$ws.Range("A21").Value = "view/metadata/custom_fields/Additional Information/Update Schedule"

# Update the selected cell shown when the sheet was last saved
$ws.Range("A22").Select()
